$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.062.32'
$ws.Range("E2").Value = '  -3.35%  '
$ws.Range("D3").Value = '3.355.73'
$ws.Range("E3").Value = '  -3.67%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '558.62'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.36'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -6.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -3.84%  '
$ws.Range("D8").Value = '3.343.08'
$ws.Range("E8").Value = '  -3.83%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.620'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.151'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.99'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("E13").Value = '  -1.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.92'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -3.59%  '
$ws.Range("D15").Value = '3.875.21'
$ws.Range("E15").Value = '  -4.42%  '
$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.118'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.44%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.353.23'
$ws.Range("E17").Value = '  -4.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.81'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.75'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -1.87%  '
$ws.Range("D20").Value = '63.983.26'
$ws.Range("E20").Value = '  -3.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.976'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '403.90'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.04%  '
$ws.Range("E23").Value = '  +0.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.27'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +3.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.33'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +9.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '82.96'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -3.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.69'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -1.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.75'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.59%  '
$ws.Range("E29").Value = '  -3.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.25'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -2.69%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.43'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '583.94'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -7.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.39'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.58%  '
$ws.Range("E34").Value = '  -3.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '58.11'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -2.91%  '
$ws.Range("E36").Value = '  -0.37%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '35.76'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.57%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.46'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("D40").Value = '0.0₃0743'
$ws.Range("E40").Value = '  -7.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.370'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.62%  '
$ws.Range("D42").Value = '3.149.06'
$ws.Range("E42").Value = '  +1.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.37%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.85'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +0.61%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.25'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.55%  '
$ws.Range("E46").Value = '  -4.18%  '
$ws.Range("E47").Value = '  -1.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.64'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -3.86%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.129'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.97%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '132.70'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.11'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.08%  '
